$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1217
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 1753
$ws.Range("L4").Value = "16 Weeks"
$ws.Range("K5").Value = 3703
$ws.Range("K6").Value = 115394
$ws.Range("K7").Value = 190758
$ws.Range("K8").Value = 15972
$ws.Range("K9").Value = 106272
$ws.Range("K10").Value = 3091
$ws.Range("K11").Value = 25485
$ws.Range("E12").Value = "JP1 JP2 JP3"
$ws.Range("H12").Value = 3
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "$0.45"
$ws.Range("J12").Style = "Normal"
$ws.Range("K12").Value = 16627
$ws.Range("K13").Value = 15549
$ws.Range("K14").Value = 335394
$ws.Range("K15").Value = 4556
$ws.Range("K16").Value = 14782
$ws.Range("K17").Value = 4392
$ws.Range("K18").Value = 3660
$ws.Range("K19").Value = 752201
$ws.Range("K20").Value = 1290726
$ws.Range("K21").Value = 1124
$ws.Range("K22").Value = 191048
$ws.Range("K23").Value = 3378
$ws.Range("K24").Value = 657
